$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B, C, D, E, F, H, J across rows 2-25 (data rows 0-23)
$data = @{
  2 = @(2.040902506846123; 0.4566974898804972; 0.03342088562028067; 0.09953734412203374; 4.747147292256642; 0.07973214163530429; 0.2477278073313158)
  3 = @(1.944238781831018; 0.4308694591487665; 0.0332331319724446; 0.09944323500848284; 4.569403720894144; 0.07973214163530429; 0.2409718732440496)
  4 = @(1.886734367927602; 0.4154811655076855; 0.03314560468927041; 0.09946038421886882; 4.462108826433024; 0.07973214163530429; 0.2369674581043597)
  5 = @(1.863761877539503; 0.409327414213692; 0.03311680360529223; 0.09948619738364073; 4.418837834462039; 0.07973214163530429; 0.2353714024204976)
  6 = @(1.85997505621549; 0.408312628618944; 0.03311243270172781; 0.09949162013109714; 4.411679755150089; 0.07973214163530429; 0.2351085286131536)
  7 = @(1.886422690812481; 0.4153977011965253; 0.03314518861003535; 0.09946065614474975; 4.461523439782894; 0.07973214163530429; 0.2369457887704414)
  8 = @(2.007187204375498; 0.4476936866682308; 0.03335031991365511; 0.09948933134743498; 4.685473094683658; 0.07973214163530429; 0.2453682961404269)
  9 = @(2.25884513608662; 0.5148137540719517; 0.03397793590968234; 0.1001412186587949; 5.139728110102169; 0.07973214163530429; 0.2630435676028782)
  10 = @(2.453063006835805; 0.566522460298188; 0.03458366768172283; 0.1009852513941496; 5.483382930874399; 0.07973214163530429; 0.2767630222231929)
  11 = @(2.543506417202082; 0.5905860372108123; 0.0348922298321952; 0.1014489930892886; 5.642031787625513; 0.07973214163530429; 0.2831696934214847)
  12 = @(2.578060842261095; 0.5997776360854346; 0.03501395315173994; 0.1016361121986229; 5.702453857129399; 0.07973214163530429; 0.2856200257044321)
  13 = @(2.570605268900181; 0.5977945128634587; 0.03498751882540319; 0.1015953002134466; 5.689425341011827; 0.07973214163530429; 0.2850912163601578)
  14 = @(2.546343081978478; 0.5913406389540796; 0.0349021455745202; 0.1014641565882215; 5.646995761260882; 0.07973214163530429; 0.2833707948191346)
  15 = @(2.531521710179845; 0.5873978185241526; 0.03485049110317817; 0.1013853274416121; 5.621051725893608; 0.07973214163530429; 0.2823201604593351)
  16 = @(2.447194812653265; 0.5649608582937162; 0.03456417780150645; 0.1009565533944112; 5.473062512912861; 0.07973214163530429; 0.276347702050856)
  17 = @(2.39600217028277; 0.5513360935781861; 0.03439708108678019; 0.1007139751532975; 5.382878081510398; 0.07973214163530429; 0.2727265280354061)
  18 = @(2.366754213258105; 0.5435503239974651; 0.03430407075791919; 0.1005819581147449; 5.331223898052514; 0.07973214163530429; 0.2706592862236903)
  19 = @(2.356885004639935; 0.540922883831513; 0.03427310742496559; 0.1005385477623157; 5.313771681972213; 0.07973214163530429; 0.2699620129090761)
  20 = @(2.401431324697967; 0.5527812015207019; 0.03441454711504832; 0.1007390207598; 5.392455771388796; 0.07973214163530429; 0.2731103943979463)
  21 = @(2.553461147161443; 0.5932341357348037; 0.03492708834155422; 0.1015023639583781; 5.659448893999695; 0.07973214163530429; 0.2838754621579227)
  22 = @(2.654604143100244; 0.6201350695391739; 0.03529055585809715; 0.1020683622280139; 5.835961989752718; 0.07973214163530429; 0.2910527067989932)
  23 = @(2.600457565000625; 0.6057347166830596; 0.0350939156487442; 0.1017601248553177; 5.741565063850715; 0.07973214163530429; 0.2872089641973332)
  24 = @(2.398976233155679; 0.5521277210328321; 0.03440664121565362; 0.10072767445617; 5.388125094268929; 0.07973214163530429; 0.2729368029132075)
  25 = @(2.189148089425828; 0.496241710876518; 0.03378328206603243; 0.09990090330498091; 5.015151403091835; 0.07973214163530429; 0.2581351263488756)
}

$cols = @("B","C","D","E","F","H","J")

foreach ($row in $data.Keys) {
  $vals = $data[$row]
  for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range("$($cols[$i])$row").Value = $vals[$i]
  }
}

Write-Output "Done updating cells"